# This script inserts two new weekly data rows (the week of 2022-04-07) at the
# top of the data table (rows 129-130) of the single data sheet, pushing the
# existing rows 129-250 down to rows 131-252. The workbook's <dimension>
# consequently grows from A1:R250 to A1:R252.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the current row 129. Excel shifts the
# formatting (including the date number format used in column D) down with
# the existing rows, so the two newly inserted rows inherit the same
# formatting that row 129 had before the insert.
$ws.Rows("129:130").Insert()

# Common values shared by every data row in this table.
$mercadoId   = 7
$mercado     = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$categoriaId = 100112023
$categoria   = "Brócoli"
$variedad    = "Sin especificar"
$unidadCom   = "`$/unidad"
$kgOUnidades = 1
$clasif      = "Hortaliza"

# --- New row 129 ---
$ws.Range("A129").Value2 = $mercadoId
$ws.Range("B129").Value  = $mercado
$ws.Range("C129").Value  = $region
$ws.Range("D129").Value2 = 44658
$ws.Range("E129").Value2 = $codreg
$ws.Range("F129").Value2 = $categoriaId
$ws.Range("G129").Value  = $categoria
$ws.Range("H129").Value  = $variedad
$ws.Range("I129").Value  = "Primera"
$ws.Range("J129").Value2 = 120
$ws.Range("K129").Value2 = 800
$ws.Range("L129").Value2 = 850
$ws.Range("M129").Value2 = 825
$ws.Range("N129").Value  = $unidadCom
$ws.Range("O129").Value  = "Provincia de Diguillín"
$ws.Range("P129").Value2 = 825
$ws.Range("Q129").Value2 = $kgOUnidades
$ws.Range("R129").Value  = $clasif

# --- New row 130 ---
$ws.Range("A130").Value2 = $mercadoId
$ws.Range("B130").Value  = $mercado
$ws.Range("C130").Value  = $region
$ws.Range("D130").Value2 = 44658
$ws.Range("E130").Value2 = $codreg
$ws.Range("F130").Value2 = $categoriaId
$ws.Range("G130").Value  = $categoria
$ws.Range("H130").Value  = $variedad
$ws.Range("I130").Value  = "Segunda"
$ws.Range("J130").Value2 = 100
$ws.Range("K130").Value2 = 700
$ws.Range("L130").Value2 = 750
$ws.Range("M130").Value2 = 725
$ws.Range("N130").Value  = $unidadCom
$ws.Range("O130").Value  = "Provincia de Diguillín"
$ws.Range("P130").Value2 = 725
$ws.Range("Q130").Value2 = $kgOUnidades
$ws.Range("R130").Value  = $clasif
